$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data table lost two rows of raw observations: "RM 232"
# (previously row 26) and "SC 92" (previously row 28). Deleting the
# entire rows shifts everything below them up automatically, which
# reproduces the row 26-33 content of the diff (and shrinks the used
# range from A1:F35 to A1:F33).
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# After the shift, "SC 119" lands on row 29 and its computed column F
# value is cleared (now blank, matching the missing-data pattern),
# while "SC 232" lands on row 33 and now has a computed column F value.
$ws.Range("F29").Value = ""
$ws.Range("F33").Value = 17.53
